# Update James Neesham's per-innings stats (runs/balls/sixes) for rows 3 and 4.
# Source sheet stores these numeric-looking values as text, so writes are
# prefixed with a leading apostrophe to keep them stored as text (matching
# the workbook's existing "number stored as text" convention) rather than
# letting Excel auto-convert them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: runs 7 -> 10, balls 7 -> 8, fours stays 0, sixes 0 -> 1
$ws.Range("C3").Value = "'10"
$ws.Range("D3").Value = "'8"
$ws.Range("F3").Value = "'1"

# Row 4: runs 10 -> 7, balls 8 -> 7, fours stays 0, sixes 1 -> 0
$ws.Range("C4").Value = "'7"
$ws.Range("D4").Value = "'7"
$ws.Range("F4").Value = "'0"
